{"js": "// Split the text \"3rd\" into three runs:\n//   \"3\"                -> normal\n//   \"rd\"                -> superscript\n//   \" im changing now\"  -> normal (appended to the end of the paragraph)\n//\n// This mirrors a user selecting \"rd\" inside \"3rd\", applying the\n// \"Superscript\" formatting (Ctrl+Shift+=), then clicking at the end of the\n// line and typing \" im changing now\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document has a single paragraph containing the text \"3rd\".\nconst paragraph = paragraphs.items[0];\n\n// Find the \"rd\" substring inside that paragraph and mark it superscript.\nconst found = paragraph.search(\"rd\", { matchCase: true, matchWholeWord: false });\nfound.load(\"items\");\nawait context.sync();\n\nconst rdRange = found.items[0];\nrdRange.font.superscript = true;\nawait context.sync();\n\n// Append the new sentence fragment at the end of the paragraph, in the\n// regular (non-superscript) style matching the rest of the paragraph's\n// language formatting.\nconst appended = paragraph.insertText(\" im changing now\", Word.InsertLocation.end);\n// Newly inserted text does not automatically inherit the run-level\n// language mark; set it explicitly to match the rest of the paragraph\n// (the BCP-47 tag, same as the document's \"en-US\" language setting).\nappended.languageId = \"en-US\";\nawait context.sync();\n", "ps1": "# Split the text \"3rd\" into three runs:\n#   \"3\"                -> normal\n#   \"rd\"                -> superscript\n#   \" im changing now\"  -> normal (appended at the end of the paragraph)\n#\n# This mirrors a user selecting \"rd\" inside \"3rd\", applying the\n# \"Superscript\" formatting (Ctrl+Shift+=), then clicking at the end of the\n# line and typing \" im changing now\".\n\n$d = $word.ActiveDocument\n\n# Find \"rd\" inside the (only) paragraph's text and make it superscript.\n$rdRange = $d.Content\n$rdRange.Find.Execute(\"rd\") | Out-Null\n$rdRange.Font.Superscript = $true\n\n# Move the cursor to the very end of the document and type the new text.\n$sel = $word.Selection\n$sel.EndKey(6)  # wdStory\n$sel.TypeText(\" im changing now\")\n\n# A freshly typed run does not automatically pick up an explicit\n# w:lang on save, so stamp it explicitly to match the rest of the\n# paragraph (same \"en-US\" language already used there).\n$newTextRange = $d.Content\n$newTextRange.Find.Execute(\" im changing now\") | Out-Null\n$newTextRange.LanguageID = \"en-US\"\n"}
